$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9390297532081604
$ws.Range("B1").Value = 1.173056244850159
$ws.Range("C1").Value = 5.240593433380127
$ws.Range("D1").Value = 1.626090288162231
$ws.Range("E1").Value = 0.95164555311203
